# Commit: "#5: fund, bonds, otherbonds, antique done"
#
# Sheet7 (基金受益憑證 / fund) gets brought up to the same shape as the
# other "done" sheets (land/building/car/stock/...):
#   - Row 1 becomes a real header row (name, owner, dealer, quantity,
#     face_value, currency, total, property_category, category, date,
#     legislator_name, legislator_id, source_file, index) instead of a
#     stray duplicate of row 2.
#   - Each data row (2-7) gains the metadata columns I:O -
#     property_category=fund, category=normal, date=2013-11-22,
#     legislator_name=吳育昇, legislator_id=1322, source_file=tmp88481,
#     index=<row id, matching column A>.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# ---- Row 1: proper header labels (B1:O1) ----------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"
$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# Match the header row's existing look (bold, centered, boxed) by copying
# the formatting already used on B1:H1 onto the newly-occupied I1:O1 cells.
$ws.Range("B1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows 2-7: add metadata columns I:O -------------------------------
$rows = @(2, 3, 4, 5, 6, 7)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "fund"                # I - property_category
    $ws.Cells.Item($r, 10).Value = "normal"              # J - category
    $ws.Cells.Item($r, 12).Value = "吳育昇"               # L - legislator_name
    $ws.Cells.Item($r, 13).Value = 1322                  # M - legislator_id
    $ws.Cells.Item($r, 14).Value = "tmp88481"            # N - source_file
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 1).Value()  # O - index (mirrors column A)

    # K - date: write it through a formula then freeze it to a literal
    # value, which keeps it a plain text cell ("2013-11-22") instead of
    # letting the smart-text parser turn it into a date serial number.
    $ws.Cells.Item($r, 11).Formula = '="2013-11-22"'
    $ws.Cells.Item($r, 11).Copy()
    $ws.Cells.Item($r, 11).PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Copy the existing data-row formatting (B2:H2) onto the new I:O cells too.
$ws.Range("B2").Copy()
$ws.Range("I2:O7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Update the worksheet dimension to cover the new columns --------------
$ws.Range("A1:O7").Select()
